# Insert 7 historical rows (2019-11-18 .. 2019-11-28) before the existing
# 2019-12-03 row, shifting all subsequent rows down by 7 (old row 129 -> new
# row 136, ..., old row 163 -> new row 170). Everything else is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows at row 129; existing rows 129:163 shift down to 136:170
# carrying their original values/formatting with them automatically.
$ws.Range("A129:A135").EntireRow.Insert()

# New row 129: 2019-11-18
$ws.Cells.Item(129, 1).Value = 1574035200
$ws.Cells.Item(129, 2).Value = "'2019-11-18"
$ws.Cells.Item(129, 3).Value = "'03023"
$ws.Cells.Item(129, 4).Value = "SMILE"
$ws.Cells.Item(129, 5).Value = 0.075
$ws.Cells.Item(129, 6).Value = 0.08
$ws.Cells.Item(129, 7).Value = 0.075
$ws.Cells.Item(129, 8).Value = 0.08
$ws.Cells.Item(129, 9).Value = 660000

# New row 130: 2019-11-19
$ws.Cells.Item(130, 1).Value = 1574121600
$ws.Cells.Item(130, 2).Value = "'2019-11-19"
$ws.Cells.Item(130, 3).Value = "'03023"
$ws.Cells.Item(130, 4).Value = "SMILE"
$ws.Cells.Item(130, 5).Value = 0.075
$ws.Cells.Item(130, 6).Value = 0.075
$ws.Cells.Item(130, 7).Value = 0.075
$ws.Cells.Item(130, 8).Value = 0.075
$ws.Cells.Item(130, 9).Value = 600000

# New row 131: 2019-11-20
$ws.Cells.Item(131, 1).Value = 1574208000
$ws.Cells.Item(131, 2).Value = "'2019-11-20"
$ws.Cells.Item(131, 3).Value = "'03023"
$ws.Cells.Item(131, 4).Value = "SMILE"
$ws.Cells.Item(131, 5).Value = 0.075
$ws.Cells.Item(131, 6).Value = 0.075
$ws.Cells.Item(131, 7).Value = 0.075
$ws.Cells.Item(131, 8).Value = 0.075
$ws.Cells.Item(131, 9).Value = 523000

# New row 132: 2019-11-21
$ws.Cells.Item(132, 1).Value = 1574294400
$ws.Cells.Item(132, 2).Value = "'2019-11-21"
$ws.Cells.Item(132, 3).Value = "'03023"
$ws.Cells.Item(132, 4).Value = "SMILE"
$ws.Cells.Item(132, 5).Value = 0.075
$ws.Cells.Item(132, 6).Value = 0.08
$ws.Cells.Item(132, 7).Value = 0.075
$ws.Cells.Item(132, 8).Value = 0.08
$ws.Cells.Item(132, 9).Value = 500000

# New row 133: 2019-11-25
$ws.Cells.Item(133, 1).Value = 1574640000
$ws.Cells.Item(133, 2).Value = "'2019-11-25"
$ws.Cells.Item(133, 3).Value = "'03023"
$ws.Cells.Item(133, 4).Value = "SMILE"
$ws.Cells.Item(133, 5).Value = 0.08
$ws.Cells.Item(133, 6).Value = 0.08
$ws.Cells.Item(133, 7).Value = 0.08
$ws.Cells.Item(133, 8).Value = 0.08
$ws.Cells.Item(133, 9).Value = 240500

# New row 134: 2019-11-26
$ws.Cells.Item(134, 1).Value = 1574726400
$ws.Cells.Item(134, 2).Value = "'2019-11-26"
$ws.Cells.Item(134, 3).Value = "'03023"
$ws.Cells.Item(134, 4).Value = "SMILE"
$ws.Cells.Item(134, 5).Value = 0.075
$ws.Cells.Item(134, 6).Value = 0.075
$ws.Cells.Item(134, 7).Value = 0.075
$ws.Cells.Item(134, 8).Value = 0.075
$ws.Cells.Item(134, 9).Value = 500000

# New row 135: 2019-11-28
$ws.Cells.Item(135, 1).Value = 1574899200
$ws.Cells.Item(135, 2).Value = "'2019-11-28"
$ws.Cells.Item(135, 3).Value = "'03023"
$ws.Cells.Item(135, 4).Value = "SMILE"
$ws.Cells.Item(135, 5).Value = 0.075
$ws.Cells.Item(135, 6).Value = 0.08
$ws.Cells.Item(135, 7).Value = 0.075
$ws.Cells.Item(135, 8).Value = 0.08
$ws.Cells.Item(135, 9).Value = 380500
